$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 43327.668
$ws.Range("J3").Value = 43327.668
$ws.Range("L3").Value = 43327.668
$ws.Range("N3").Value = -43555.668

$ws.Range("H17").Value = 41667744
$ws.Range("J17").Value = 1066.6818
$ws.Range("L17").Value = 3200.0454
$ws.Range("N17").Value = -3536.0454

$ws.Range("H76").Value = 4475
$ws.Range("J76").Value = 4400
$ws.Range("L76").Value = 4400
$ws.Range("N76").Value = -5030

$ws.Range("H79").Value = 4475
$ws.Range("J79").Value = 4400
$ws.Range("L79").Value = 4400
$ws.Range("N79").Value = -6584

$ws.Range("H93").Value = 69000
$ws.Range("J93").Value = 69000
$ws.Range("L93").Value = 69000
$ws.Range("N93").Value = -73992

$ws.Range("H102").Value = 43327.668
$ws.Range("J102").Value = 43327.668
$ws.Range("L102").Value = 43327.668
$ws.Range("N102").Value = -49817.668

$ws.Range("H130").Value = 19997.143
$ws.Range("J130").Value = 19997.143
$ws.Range("L130").Value = 19997.143
$ws.Range("N130").Value = -30037.143

$ws.Range("H137").Value = 11125926
$ws.Range("I137").Value = 20025024
$ws.Range("K137").Value = 60075072
$ws.Range("M137").Value = -60072522

$ws.Range("H138").Value = 4702.3423
$ws.Range("I138").Value = 3124.1667
$ws.Range("J138").Value = 4998.25
$ws.Range("K138").Value = 9372.500100000001
$ws.Range("L138").Value = 14994.75
$ws.Range("M138").Value = -4232.500100000001
$ws.Range("N138").Value = -25274.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14497077
$ws.Range("I32").Value = 15877618
$ws.Range("J32").Value = 1402.1666
$ws.Range("K32").Value = 15877618
$ws.Range("L32").Value = 1402.1666
$ws.Range("M32").Value = -15877331
$ws.Range("N32").Value = -1976.1666

$ws.Range("H132").Value = 1238.56
$ws.Range("I132").Value = 1043.2174
$ws.Range("K132").Value = 3129.6522
$ws.Range("M132").Value = -599.6522

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 54497.156
$ws.Range("I20").Value = 1877.3636
$ws.Range("K20").Value = 1877.3636
$ws.Range("M20").Value = -1630.3636

$ws.Range("H59").Value = 77500
$ws.Range("J59").Value = 77500
$ws.Range("L59").Value = 77500
$ws.Range("N59").Value = -79194

$ws.Range("H99").Value = 5913.625
$ws.Range("I99").Value = 2261.8
$ws.Range("K99").Value = 2261.8
$ws.Range("M99").Value = -763.8000000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5246.7744
$ws.Range("I31").Value = 7385.9375
$ws.Range("K31").Value = 7385.9375
$ws.Range("M31").Value = -7090.9375

$ws.Range("H32").Value = 1047
$ws.Range("I32").Value = 1047
$ws.Range("K32").Value = 1047
$ws.Range("M32").Value = -731

$ws.Range("H34").Value = 5246.7744
$ws.Range("I34").Value = 7385.9375
$ws.Range("K34").Value = 7385.9375
$ws.Range("M34").Value = -7183.9375

$ws.Range("H47").Value = 36690.332
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 36690.332
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 36690.332
$ws.Range("N47").Value = -37822.332
$ws.Range("M47").ClearContents()

$ws.Range("H94").Value = 1711.6154
$ws.Range("I94").Value = 1171.3334
$ws.Range("J94").Value = 1873.7
$ws.Range("K94").Value = 1171.3334
$ws.Range("L94").Value = 1873.7
$ws.Range("M94").Value = -720.3334
$ws.Range("N94").Value = -2775.7

$ws.Range("H99").Value = 24934192
$ws.Range("I99").Value = 3486531.8
$ws.Range("K99").Value = 3486531.8
$ws.Range("M99").Value = -3485033.8

$ws.Range("H126").Value = 24934192
$ws.Range("I126").Value = 3486531.8
$ws.Range("K126").Value = 10459595.4
$ws.Range("M126").Value = -10457125.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 40.055557
$ws.Range("I2").Value = 19.153847
$ws.Range("K2").Value = 114.923082
$ws.Range("M2").Value = -1.923081999999994

$ws.Range("H17").Value = 193.88889
$ws.Range("I17").Value = 128.75
$ws.Range("J17").Value = 246
$ws.Range("K17").Value = 386.25
$ws.Range("L17").Value = 738
$ws.Range("M17").Value = -217.25
$ws.Range("N17").Value = -1076

$ws.Range("H55").Value = 15608918
$ws.Range("J55").Value = 18073258
$ws.Range("L55").Value = 54219774
$ws.Range("N55").Value = -54220128

$ws.Range("H131").Value = 1430949.1
$ws.Range("I131").Value = 2000718.8
$ws.Range("K131").Value = 6002156.4
$ws.Range("M131").Value = -5997116.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4920.25
$ws.Range("I126").Value = 4566.75
$ws.Range("J126").Value = 5273.75
$ws.Range("K126").Value = 13700.25
$ws.Range("L126").Value = 15821.25
$ws.Range("M126").Value = -11230.25
$ws.Range("N126").Value = -20761.25

$ws.Range("H132").Value = 6621.4326
$ws.Range("I132").Value = 6193.5483
$ws.Range("K132").Value = 18580.6449
$ws.Range("M132").Value = -16050.6449

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2347.303
$ws.Range("J22").Value = 1776
$ws.Range("L22").Value = 1776
$ws.Range("N22").Value = -2366

$ws.Range("H27").Value = 2347.303
$ws.Range("J27").Value = 1776
$ws.Range("L27").Value = 1776
$ws.Range("N27").Value = -1990

$ws.Range("H93").Value = 1487.15
$ws.Range("I93").Value = 1491.7646
$ws.Range("K93").Value = 1491.7646
$ws.Range("M93").Value = -243.7646

$ws.Range("H100").Value = 74547.94500000001
$ws.Range("I100").Value = 94635.62
$ws.Range("K100").Value = 94635.62
$ws.Range("M100").Value = -94094.62

$ws.Range("H122").Value = 4770.643
$ws.Range("I122").Value = 4390
$ws.Range("K122").Value = 13170
$ws.Range("M122").Value = -10720

$ws.Range("H132").Value = 4776.087
$ws.Range("I132").Value = 4954.425
$ws.Range("J132").Value = 3587.1667
$ws.Range("K132").Value = 14863.275
$ws.Range("L132").Value = 10761.5001
$ws.Range("M132").Value = -12333.275
$ws.Range("N132").Value = -15821.5001

$ws.Range("H136").Value = 5711.4165
$ws.Range("I136").Value = 15645.667
$ws.Range("K136").Value = 46937.001
$ws.Range("M136").Value = -44387.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 15014
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H93").Value = 29994
$ws.Range("J93").Value = 29994
$ws.Range("L93").Value = 29994
$ws.Range("N93").Value = -34986

$ws.Range("H103").Value = 57883.6
$ws.Range("J103").Value = 62354.5
$ws.Range("L103").Value = 62354.5
$ws.Range("N103").Value = -64698.5

$ws.Range("H132").Value = 4093.2114
$ws.Range("I132").Value = 3463.4634
$ws.Range("K132").Value = 10390.3902
$ws.Range("M132").Value = -7860.3902

$ws.Range("H136").Value = 1617370.5
$ws.Range("I136").Value = 2503068.5
$ws.Range("J136").Value = 7010.636
$ws.Range("K136").Value = 7509205.5
$ws.Range("L136").Value = 21031.908
$ws.Range("M136").Value = -7506655.5
$ws.Range("N136").Value = -26131.908

$ws.Range("H138").Value = 75999.8
$ws.Range("J138").Value = 75999.8
$ws.Range("L138").Value = 75999.8
$ws.Range("N138").Value = -86279.8
